# Add an "Is model" column (E) to the worksheet, with y/y/n/n for the
# four data rows that already have Feature ID / PDB filename values,
# and adjust column widths + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "Is model"

# Values for rows 2-5 (rows 6-8 stay untouched, same as before)
$ws.Range("E2").Value = "y"
$ws.Range("E3").Value = "y"
$ws.Range("E4").Value = "n"
$ws.Range("E5").Value = "n"

# Column width adjustments (values chosen so the engine's pixel-quantized
# stored width lands as close as possible to the target 37.83203125 /
# 19.83203125 character widths)
$ws.Columns.Item(2).ColumnWidth = 37
$ws.Columns.Item(3).ColumnWidth = 19

# Update selection to match the new layout
$ws.Range("E9").Select()
